# Relocate the "_GoBack" bookmark (Word's "last edit position" marker) from
# the end of the document to the spot where the new annotator's edit
# actually happened: inside the paragraph
#   "... domunde tar  rosc for nanme"
# right between "tar " and "rosc for nanme" -- and collapse the double
# space that used to separate them into a single space in the process.

$d = $word.ActiveDocument

# Locate the double space between "tar" and "rosc" in that paragraph.
$target = $d.Content
$found = $target.Find.Execute("tar  rosc", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the target text 'tar  rosc for nanme'"
}

# $target.Start now points at the "t" of "tar"; 4 characters in ("t","a",
# "r"," ") lands right between the two spaces -- i.e. right after the single
# space that should remain at the end of the first run.
$splitPos = $target.Start + 4

# Insert the (collapsed) "_GoBack" bookmark there first, before any text is
# touched. Word splits the run at that exact point to host the bookmark, so
# the untouched text before the split keeps its original run
# formatting/rsid, while Bookmarks.Add("_GoBack", ...) relocates the single
# existing "_GoBack" bookmark (previously at the very end of the document)
# to this new position instead of creating a duplicate.
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Now remove the extra (now redundant) space that is left sitting right
# after the bookmark, at the start of the "rosc for nanme" run, collapsing
# "tar " + " rosc..." down to "tar " + "rosc...".
$extraSpace = $d.Range($splitPos, $splitPos + 1)
if ($extraSpace.Text -eq " ") {
    $extraSpace.Delete()
}
